# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AC1=Wins, AD1=Losses, AE1=Ties ---
# Copy the existing header style (bold, centered, bordered) from AB1
# so the new header cells share the same style index as the rest of
# the header row, then set their text.
$ws.Range("AB1").Copy($ws.Range("AC1"))
$ws.Range("AB1").Copy($ws.Range("AD1"))
$ws.Range("AB1").Copy($ws.Range("AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2-37): every player on this roster shares the same
# 1993 team record: 81 wins, 81 losses, 0 ties ---
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 29).Value = 81   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 81   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
